$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 3
$ws.Range("B7").Value = 4
$ws.Range("B22").Value = 3
$ws.Range("B84").Value = 3
